$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.164.95'
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").Value = '2.211.24'
$ws.Range("E3").Value = '  -0.51%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = "'295.41"
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("D6").Value = "'87.75"
$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("D7").Value = "'0.513"
$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -0.50%  '

$ws.Range("D10").Value = "'30.87"
$ws.Range("E10").Value = '  +1.68%  '

$ws.Range("D11").Value = "'51.43"
$ws.Range("E11").Value = '  +5.79%  '

$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("E13").Value = '  +2.45%  '

$ws.Range("E14").Value = '  -1.57%  '

$ws.Range("D15").Value = '2.552.65'
$ws.Range("E15").Value = '  -0.66%  '

$ws.Range("D16").Value = "'13.81"
$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("D17").Value = '2.159.02'
$ws.Range("E17").Value = '  -2.51%  '

$ws.Range("E18").Value = '  +1.21%  '

$ws.Range("D19").Value = '40.069.51'
$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("E20").Value = '  +0.19%  '

$ws.Range("D21").Value = "'11.30"
$ws.Range("E21").Value = '  -0.99%  '

$ws.Range("E22").Value = '  -1.19%  '

$ws.Range("D23").Value = "'65.59"
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").Value = "'235.46"

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").Value = "'2.48"
$ws.Range("E26").Value = '  +0.68%  '

$ws.Range("E27").Value = '  -0.70%  '

$ws.Range("D28").Value = "'23.20"
$ws.Range("E28").Value = '  +2.26%  '

$ws.Range("D29").Value = "'9.33"
$ws.Range("E29").Value = '  +1.30%  '

$ws.Range("E30").Value = '  -4.81%  '

$ws.Range("D31").Value = "'156.44"
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("D32").Value = "'32.06"
$ws.Range("E32").Value = '  +0.96%  '

$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("D34").Value = "'4.96"
$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("D35").Value = "'3.01"
$ws.Range("E35").Value = '  +4.69%  '

$ws.Range("E36").Value = '  -0.52%  '

$ws.Range("E37").Value = '  -0.76%  '

$ws.Range("E38").Value = '  +1.53%  '

$ws.Range("E39").Value = '  +3.04%  '

$ws.Range("E40").Value = '  +2.48%  '

$ws.Range("D41").Value = "'15.70"
$ws.Range("E41").Value = '  +0.37%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = "'3.81"
$ws.Range("E42").Value = '  -1.22%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.077.59'
$ws.Range("E43").Value = '  -1.92%  '

$ws.Range("D44").Value = "'19.25"
$ws.Range("E44").Value = '  +7.55%  '

$ws.Range("E45").Value = '  +1.34%  '

$ws.Range("D46").Value = "'9.97"
$ws.Range("E46").Value = '  +0.57%  '

$ws.Range("E47").Value = '  +5.60%  '

$ws.Range("D49").Value = '2.426.36'
$ws.Range("E49").Value = '  -0.36%  '

$ws.Range("D50").Value = "'1.13"
$ws.Range("E50").Value = '  +2.15%  '

$ws.Range("E51").Value = '  +1.17%  '

# Reset style for text-forced numeric-looking cells to avoid quote-prefix styling artifacts
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D50").Style = "Normal"
